# Update countries & provincias Spain
# Applies the 24-Sep-2020 data refresh to the "Pais" sheet:
#  - Updates case counts for several countries (rows whose totals moved).
#  - "Malasia" overtakes "Consejo Danes para los Refugiados" in the ranking,
#    so the two rows swap positions (row 97 <-> row 98).
#  - "Islas Malvinas" and "Montserrat" swap positions too (row 214 <-> row 215).
#  - Updates the "last updated" timestamp string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: Estados Unidos ---
$ws.Range("B4").Value = 7140137
$ws.Range("C4").Value = 584
$ws.Range("D4").Value = 4399996
$ws.Range("E4").Value = 2533543
$ws.Range("G4").Value = 5
$ws.Range("H4").Value = 206598

# --- Row 5: India ---
$ws.Range("B5").Value = 5737197
$ws.Range("C5").Value = 7013
$ws.Range("E5").Value = 971006
$ws.Range("G5").Value = 31
$ws.Range("H5").Value = 91204

# --- Row 18: Banglades ---
$ws.Range("B18").Value = 355384
$ws.Range("C18").Value = 1540
$ws.Range("D18").Value = 265092
$ws.Range("E18").Value = 85220
$ws.Range("G18").Value = 28
$ws.Range("H18").Value = 5072

# --- Row 33: Rumania ---
$ws.Range("B33").Value = 118054
$ws.Range("C33").Value = 1639
$ws.Range("D33").Value = 94877
$ws.Range("E33").Value = 18586
$ws.Range("G33").Value = 41
$ws.Range("H33").Value = 4591

# --- Row 44: Emiratos Arabes Unidos ---
$ws.Range("B44").Value = 88532
$ws.Range("C44").Value = 1002
$ws.Range("D44").Value = 77937
$ws.Range("E44").Value = 10188
$ws.Range("G44").Value = 1
$ws.Range("H44").Value = 407

# --- Row 47: Polonia ---
$ws.Range("B47").Value = 82809
$ws.Range("C47").Value = 1136
$ws.Range("D47").Value = 66158
$ws.Range("E47").Value = 14282
$ws.Range("G47").Value = 25
$ws.Range("H47").Value = 2369

# --- Row 90: Senegal ---
$ws.Range("B90").Value = 14816
$ws.Range("C90").Value = 21
$ws.Range("D90").Value = 11818
$ws.Range("E90").Value = 2694
$ws.Range("G90").Value = 1
$ws.Range("H90").Value = 304

# --- Rows 97/98: Malasia overtakes Consejo Danes para los Refugiados ---
$ws.Range("A97").Value = "Malasia"
$ws.Range("B97").Value = 10576
$ws.Range("C97").Value = 71
$ws.Range("D97").Value = 9666
$ws.Range("E97").Value = 777
$ws.Range("H97").Value = 133

$ws.Range("A98").Value = "Consejo Danes para los Refugiados"
$ws.Range("B98").Value = 10537
$ws.Range("C98").Value = 0
$ws.Range("D98").Value = 10041
$ws.Range("E98").Value = 225
$ws.Range("H98").Value = 271

# --- Row 122: Hong Kong ---
$ws.Range("B122").Value = 5057
$ws.Range("C122").Value = 7
$ws.Range("E122").Value = 204
$ws.Range("G122").Value = 1
$ws.Range("H122").Value = 104

# --- Row 130: Georgia ---
$ws.Range("E130").Value = 2668
$ws.Range("G130").Value = 1
$ws.Range("H130").Value = 26

# --- Row 137: Gambia ---
$ws.Range("B137").Value = 3552
$ws.Range("C137").Value = 10
$ws.Range("D137").Value = 2012
$ws.Range("E137").Value = 1430

# --- Row 183: Gibraltar ---
$ws.Range("B183").Value = 361
$ws.Range("C183").Value = 4
$ws.Range("D183").Value = 331
$ws.Range("E183").Value = 30

# --- Rows 214/215: Islas Malvinas swaps with Montserrat ---
$ws.Range("A214").Value = "Islas Malvinas"
$ws.Range("D214").Value = 13
$ws.Range("H214").Value = 0

$ws.Range("A215").Value = "Montserrat"
$ws.Range("D215").Value = 12
$ws.Range("H215").Value = 1

# --- Timestamp update (header cell, row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 24 de Septiembre de 2020 a las 12:31"
